$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 297; this shifts the existing rows
# 297-363 down to 299-365 and carries formatting from the row below.
$ws.Rows("297:298").Insert()

# Populate the first new row (297)
$ws.Range("A297").Value = 10
$ws.Range("B297").Value = "Vega Modelo de Temuco"
$ws.Range("C297").Value = "La Araucanía"
$ws.Range("D297").Value = 44511
$ws.Range("E297").Value = 9
$ws.Range("F297").Value = "Fruta"
$ws.Range("G297").Value = 100101
$ws.Range("H297").Value = "Berries"
$ws.Range("I297").Value = 100101007
$ws.Range("J297").Value = "Kiwi"
$ws.Range("K297").Value = "Hayward"
$ws.Range("L297").Value = "Especial"
$ws.Range("M297").Value = 55
$ws.Range("N297").Value = 24000
$ws.Range("O297").Value = 24000
$ws.Range("P297").Value = 24000
$ws.Range("Q297").Value = "`$/bandeja 18 kilos"
$ws.Range("R297").Value = "Región de O'Higgins"
$ws.Range("S297").Value = 1333
$ws.Range("T297").Value = 18

# Populate the second new row (298)
$ws.Range("A298").Value = 10
$ws.Range("B298").Value = "Vega Modelo de Temuco"
$ws.Range("C298").Value = "La Araucanía"
$ws.Range("D298").Value = 44511
$ws.Range("E298").Value = 9
$ws.Range("F298").Value = "Fruta"
$ws.Range("G298").Value = 100101
$ws.Range("H298").Value = "Berries"
$ws.Range("I298").Value = 100101007
$ws.Range("J298").Value = "Kiwi"
$ws.Range("K298").Value = "Hayward"
$ws.Range("L298").Value = "Primera"
$ws.Range("M298").Value = 65
$ws.Range("N298").Value = 15000
$ws.Range("O298").Value = 15000
$ws.Range("P298").Value = 15000
$ws.Range("Q298").Value = "`$/bandeja 10 kilos"
$ws.Range("R298").Value = "Región de O'Higgins"
$ws.Range("S298").Value = 1500
$ws.Range("T298").Value = 10
